$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row in column A (data currently spans A1:A13)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2025-09-08 10:27:25"
